$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new bug report row (row 6) with the same structure as the existing rows.
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Solve command has no input for empty locations"
$ws.Cells.Item(6, 3).Value = "ricky"
$ws.Cells.Item(6, 4).Value = "No else statement to handle"
$ws.Cells.Item(6, 5).Value = "Else statement to handle empty locations added"
$ws.Cells.Item(6, 6).Value = "Fixed"

# Update the active selection to match the target state.
$ws.Range("F6").Select()
